$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": budget/sales figures for two clients move from 0
# to a real amount, and the "N de 57" progress counters for columns C/D bump
# up by one (2->3, 1->2) reflecting the newly counted entries.
# ---------------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("C29").Value = 388.8
$wsGrupo.Range("D33").Value = 1068.48
$wsGrupo.Range("C59").Value = "3 de 57"
$wsGrupo.Range("D59").Value = "2 de 57"

# ---------------------------------------------------------------------------
# Sheet "VENTA MENSUAL": same two clients' monthly sales (column F) move from
# 0 to the matching amount, and the column total (F59) increases accordingly.
# ---------------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F29").Value = 388.8
$wsMensual.Range("F33").Value = 1068.48
$wsMensual.Range("F59").Value = 6195.34

# ---------------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL": VENTA (D), POR CUMPLIR (E) and CUMPLIMIENTO
# (F) are recomputed for the affected groups (rows 2, 3) and the totals
# row (14), plus the PORCELANATO group subtotal (row 12) which rolls up the
# per-group figures.
#   POR CUMPLIR   = PRESUPUESTO - VENTA   (E = C - D)
#   CUMPLIMIENTO  = VENTA / PRESUPUESTO   (F = D / C)
# ---------------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumpl.Range("D2").Value = 777.6
$wsCumpl.Range("E2").Value = 2114.60588040374
$wsCumpl.Range("F2").Value = 0.2688605279688631

$wsCumpl.Range("D3").Value = 1353.41
$wsCumpl.Range("E3").Value = 19034.0674217135
$wsCumpl.Range("F3").Value = 0.06638437762577547

$wsCumpl.Range("D12").Value = 4660.88
$wsCumpl.Range("E12").Value = 43963.18
$wsCumpl.Range("F12").Value = 0.09585542630541342

$wsCumpl.Range("D14").Value = 8601.74
$wsCumpl.Range("E14").Value = 91296.25284188786
$wsCumpl.Range("F14").Value = 0.08610523350168089

# The updated (shorter) "POR CUMPLIR" values let column E's auto-fit width
# shrink from 23 to 22 characters.
$wsCumpl.Columns.Item(5).ColumnWidth = 21.166666666666668
